$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.189171195030212
$ws.Range("B1").Value = 1.934024691581726
$ws.Range("C1").Value = 6.454429149627686
$ws.Range("D1").Value = 2.288463115692139
$ws.Range("E1").Value = 1.194097876548767
